# Regenerate the "K" column (G) values for save_data sheet.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the header is already "K"; this updates the
# per-row numeric K values that were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newValues = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 3
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
